$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5469.8237
$ws.Range("I76").Value = 5499.1333
$ws.Range("K76").Value = 5499.1333
$ws.Range("M76").Value = -5184.1333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5469.8237
$ws.Range("I79").Value = 5499.1333
$ws.Range("K79").Value = 5499.1333
$ws.Range("M79").Value = -4407.1333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 26295.8
$ws.Range("J95").Value = 26295.8
$ws.Range("L95").Value = 26295.8
$ws.Range("N95").Value = -31787.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 15402.143
$ws.Range("I96").Value = 19331.092
$ws.Range("J96").Value = 996
$ws.Range("K96").Value = 57993.276
$ws.Range("L96").Value = 2988
$ws.Range("M96").Value = -56620.276
$ws.Range("N96").Value = -5734

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 8400
$ws.Range("I99").Value = 8400
$ws.Range("K99").Value = 25200
$ws.Range("M99").Value = -23702

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1062.8572
$ws.Range("J112").Value = 990
$ws.Range("L112").Value = 2970
$ws.Range("N112").Value = -5186

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3559
$ws.Range("I113").Value = 4156.857
$ws.Range("J113").Value = 3035.875
$ws.Range("K113").Value = 4156.857
$ws.Range("L113").Value = 3035.875
$ws.Range("M113").Value = -902.857
$ws.Range("N113").Value = -9543.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3340.6418
$ws.Range("I138").Value = 3337
$ws.Range("J138").Value = 3342.1914
$ws.Range("K138").Value = 10011
$ws.Range("L138").Value = 10026.5742
$ws.Range("M138").Value = -4871
$ws.Range("N138").Value = -20306.5742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4316.1816
$ws.Range("I32").Value = 3820.484
$ws.Range("J32").Value = 11999.5
$ws.Range("K32").Value = 3820.484
$ws.Range("L32").Value = 11999.5
$ws.Range("M32").Value = -3533.484
$ws.Range("N32").Value = -12573.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1342.9584
$ws.Range("I88").Value = 1023.6
$ws.Range("J88").Value = 1427
$ws.Range("K88").Value = 1023.6
$ws.Range("L88").Value = 1427
$ws.Range("M88").Value = -617.6
$ws.Range("N88").Value = -2239

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1342.9584
$ws.Range("I91").Value = 1023.6
$ws.Range("J91").Value = 1427
$ws.Range("K91").Value = 1023.6
$ws.Range("L91").Value = 1427
$ws.Range("M91").Value = 380.4
$ws.Range("N91").Value = -4235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1027.9412
$ws.Range("I102").Value = 1027.9412
$ws.Range("K102").Value = 1027.9412
$ws.Range("M102").Value = 594.0588

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1865.6154
$ws.Range("I122").Value = 1485.3
$ws.Range("K122").Value = 4455.9
$ws.Range("M122").Value = -2005.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 715.5714
$ws.Range("I94").Value = 762
$ws.Range("J94").Value = 599.5
$ws.Range("K94").Value = 762
$ws.Range("L94").Value = 599.5
$ws.Range("M94").Value = -311
$ws.Range("N94").Value = -1501.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 10994.5
$ws.Range("I96").Value = 10994.5
$ws.Range("K96").Value = 10994.5
$ws.Range("M96").Value = -8248.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4398.4
$ws.Range("I107").Value = 4325.4116
$ws.Range("K107").Value = 4325.4116
$ws.Range("M107").Value = -2405.4116

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6283.1763
$ws.Range("I31").Value = 5070
$ws.Range("J31").Value = 6656.4614
$ws.Range("K31").Value = 5070
$ws.Range("L31").Value = 6656.4614
$ws.Range("M31").Value = -4775
$ws.Range("N31").Value = -7246.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6283.1763
$ws.Range("I34").Value = 5070
$ws.Range("J34").Value = 6656.4614
$ws.Range("K34").Value = 5070
$ws.Range("L34").Value = 6656.4614
$ws.Range("M34").Value = -4868
$ws.Range("N34").Value = -7060.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1363.5
$ws.Range("I62").Value = 1363.5
$ws.Range("K62").Value = 1363.5
$ws.Range("M62").Value = -739.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 1363.5
$ws.Range("I65").Value = 1363.5
$ws.Range("K65").Value = 6817.5
$ws.Range("M65").Value = -3697.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11002.667
$ws.Range("I86").Value = 20000
$ws.Range("K86").Value = 20000
$ws.Range("M86").Value = -18877

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 11002.667
$ws.Range("I89").Value = 20000
$ws.Range("K89").Value = 100000
$ws.Range("M89").Value = -94384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 99994
$ws.Range("J100").Value = 99994
$ws.Range("L100").Value = 99994
$ws.Range("N100").Value = -102158

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3033.4814
$ws.Range("I105").Value = 2224.5833
$ws.Range("K105").Value = 2224.5833
$ws.Range("M105").Value = -477.5832999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 221191.6
$ws.Range("I2").Value = 550005
$ws.Range("J2").Value = 138988.25
$ws.Range("K2").Value = 3300030
$ws.Range("L2").Value = 833929.5
$ws.Range("M2").Value = -3299917
$ws.Range("N2").Value = -834155.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 489.8
$ws.Range("I60").Value = 489.8
$ws.Range("K60").Value = 1469.4
$ws.Range("M60").Value = -1218.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1117.0834
$ws.Range("I113").Value = 1065.7142
$ws.Range("J113").Value = 1138.2354
$ws.Range("K113").Value = 3197.1426
$ws.Range("L113").Value = 3414.7062
$ws.Range("M113").Value = -1027.1426
$ws.Range("N113").Value = -7754.706200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 12698.0625
$ws.Range("I121").Value = 25001.8
$ws.Range("J121").Value = 7105.4546
$ws.Range("K121").Value = 75005.39999999999
$ws.Range("L121").Value = 21316.3638
$ws.Range("M121").Value = -73695.39999999999
$ws.Range("N121").Value = -23936.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 688
$ws.Range("J122").Value = 657.25
$ws.Range("L122").Value = 5915.25
$ws.Range("N122").Value = -10815.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11878.583
$ws.Range("I80").Value = 2996
$ws.Range("J80").Value = 12686.091
$ws.Range("K80").Value = 2996
$ws.Range("L80").Value = 12686.091
$ws.Range("M80").Value = -1998
$ws.Range("N80").Value = -14682.091

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 11878.583
$ws.Range("I83").Value = 2996
$ws.Range("J83").Value = 12686.091
$ws.Range("K83").Value = 14980
$ws.Range("L83").Value = 63430.455
$ws.Range("M83").Value = -9988
$ws.Range("N83").Value = -73414.455

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 72246.25
$ws.Range("J140").Value = 77995
$ws.Range("L140").Value = 77995
$ws.Range("N140").Value = -88355

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6746.148
$ws.Range("I7").Value = 5046.231
$ws.Range("J7").Value = 8324.643
$ws.Range("K7").Value = 5046.231
$ws.Range("L7").Value = 8324.643
$ws.Range("M7").Value = -4934.231
$ws.Range("N7").Value = -8548.643

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1257.6
$ws.Range("I82").Value = 1213.1666
$ws.Range("J82").Value = 1324.25
$ws.Range("K82").Value = 1213.1666
$ws.Range("L82").Value = 1324.25
$ws.Range("M82").Value = -852.1666
$ws.Range("N82").Value = -2046.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1257.6
$ws.Range("I85").Value = 1213.1666
$ws.Range("J85").Value = 1324.25
$ws.Range("K85").Value = 1213.1666
$ws.Range("L85").Value = 1324.25
$ws.Range("M85").Value = 34.83339999999998
$ws.Range("N85").Value = -3820.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6746.148
$ws.Range("I126").Value = 5046.231
$ws.Range("J126").Value = 8324.643
$ws.Range("K126").Value = 15138.693
$ws.Range("L126").Value = 24973.929
$ws.Range("M126").Value = -12668.693
$ws.Range("N126").Value = -29913.929

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4376.273
$ws.Range("I126").Value = 2556.5715
$ws.Range("J126").Value = 7560.75
$ws.Range("K126").Value = 7669.7145
$ws.Range("L126").Value = 22682.25
$ws.Range("M126").Value = -5199.7145
$ws.Range("N126").Value = -27622.25
